$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update PLC2's IP address from 192.168.0.130 to 192.168.53.82
$ws.Range("B3").Value = "192.168.53.82"
